$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Date column (BF) for all data rows 2-31: "6-18-2014-15" -> "2015-06-18"
# Force text format so the date-like string is not auto-converted to a date serial number.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.NumberFormat = "@"
    $cell.Value = "2015-06-18"
}

# Apply corrected numeric values (training data fix, one day off)
$ws.Range("AS2").Value = 22
$ws.Range("AH3").Value = 8
$ws.Range("AQ3").Value = 14
$ws.Range("AY3").Value = 21
$ws.Range("AZ3").Value = 20
$ws.Range("AH4").Value = 3
$ws.Range("AT4").Value = 22
$ws.Range("AU4").Value = 20
$ws.Range("AR5").Value = 26
$ws.Range("AY5").Value = 23
$ws.Range("AF6").Value = 9
$ws.Range("AG6").Value = 9
$ws.Range("AH6").Value = 8
$ws.Range("AV6").Value = 12
$ws.Range("AH7").Value = 23
$ws.Range("AP7").Value = 12
$ws.Range("AQ7").Value = 18
$ws.Range("AS7").Value = 21
$ws.Range("D8").Value = 82
$ws.Range("F8").Value = 32
$ws.Range("G8").Value = 0.61
$ws.Range("J8").Value = 85.8
$ws.Range("K8").Value = 0.463
$ws.Range("N8").Value = 0.352
$ws.Range("Q8").Value = 0.752
$ws.Range("S8").Value = 31.8
$ws.Range("T8").Value = 42.3
$ws.Range("U8").Value = 22.5
$ws.Range("AA8").Value = 22.1
$ws.Range("AB8").Value = 105.2
$ws.Range("AC8").Value = 2.9
$ws.Range("AD8").Value = 1
$ws.Range("AF8").Value = 9
$ws.Range("AJ8").Value = 8
$ws.Range("AQ8").Value = 16
$ws.Range("AS8").Value = 23
$ws.Range("AT8").Value = 23
$ws.Range("AX8").Value = 20
$ws.Range("AH9").Value = 8
$ws.Range("AS9").Value = 12
$ws.Range("AX9").Value = 21
$ws.Range("AJ10").Value = 6
$ws.Range("BA10").Value = 24
$ws.Range("BB10").Value = 18
$ws.Range("D12").Value = 82
$ws.Range("E12").Value = 56
$ws.Range("G12").Value = 0.6830000000000001
$ws.Range("J12").Value = 83.3
$ws.Range("N12").Value = 0.348
$ws.Range("O12").Value = 18.6
$ws.Range("P12").Value = 26
$ws.Range("Q12").Value = 0.715
$ws.Range("T12").Value = 43.7
$ws.Range("Y12").Value = 5.3
$ws.Range("Z12").Value = 22
$ws.Range("AA12").Value = 21.1
$ws.Range("AC12").Value = 3.4
$ws.Range("AD12").Value = 1
$ws.Range("AE12").Value = 3
$ws.Range("AG12").Value = 3
$ws.Range("AO12").Value = 5
$ws.Range("AP12").Value = 2
$ws.Range("AT12").Value = 14
$ws.Range("AY12").Value = 22
$ws.Range("AQ13").Value = 13
$ws.Range("D15").Value = 82
$ws.Range("F15").Value = 61
$ws.Range("G15").Value = 0.256
$ws.Range("I15").Value = 37.2
$ws.Range("J15").Value = 85.59999999999999
$ws.Range("K15").Value = 0.435
$ws.Range("M15").Value = 18.9
$ws.Range("O15").Value = 17.5
$ws.Range("P15").Value = 23.6
$ws.Range("Q15").Value = 0.741
$ws.Range("S15").Value = 32.3
$ws.Range("T15").Value = 43.9
$ws.Range("U15").Value = 20.9
$ws.Range("Y15").Value = 4.8
$ws.Range("Z15").Value = 21.2
$ws.Range("AA15").Value = 19.4
$ws.Range("AB15").Value = 98.5
$ws.Range("AC15").Value = -6.8
$ws.Range("AD15").Value = 1
$ws.Range("AP15").Value = 11
$ws.Range("AS15").Value = 13
$ws.Range("AT15").Value = 12
$ws.Range("AU15").Value = 21
$ws.Range("AX15").Value = 22
$ws.Range("AZ15").Value = 21
$ws.Range("BA15").Value = 23
$ws.Range("BB15").Value = 19
$ws.Range("AE16").Value = 5
$ws.Range("AK16").Value = 9
$ws.Range("AL17").Value = 21
$ws.Range("AX17").Value = 18
$ws.Range("AH18").Value = 3
$ws.Range("AH19").Value = 20
$ws.Range("AP19").Value = 3
$ws.Range("D20").Value = 82
$ws.Range("E20").Value = 45
$ws.Range("G20").Value = 0.549
$ws.Range("J20").Value = 82.90000000000001
$ws.Range("K20").Value = 0.457
$ws.Range("L20").Value = 7.1
$ws.Range("N20").Value = 0.37
$ws.Range("P20").Value = 21.8
$ws.Range("Q20").Value = 0.751
$ws.Range("R20").Value = 11.5
$ws.Range("S20").Value = 32
$ws.Range("T20").Value = 43.5
$ws.Range("X20").Value = 6.2
$ws.Range("AA20").Value = 18.7
$ws.Range("AB20").Value = 99.40000000000001
$ws.Range("AC20").Value = 0.8
$ws.Range("AD20").Value = 1
$ws.Range("AE20").Value = 13
$ws.Range("AG20").Value = 13
$ws.Range("AI20").Value = 11
$ws.Range("AJ20").Value = 19
$ws.Range("AK20").Value = 10
$ws.Range("AQ20").Value = 17
$ws.Range("AS20").Value = 19
$ws.Range("AW20").Value = 25
$ws.Range("BC20").Value = 13
$ws.Range("AO22").Value = 6
$ws.Range("AQ22").Value = 15
$ws.Range("D23").Value = 82
$ws.Range("F23").Value = 57
$ws.Range("G23").Value = 0.305
$ws.Range("I23").Value = 37.5
$ws.Range("K23").Value = 0.453
$ws.Range("M23").Value = 19.5
$ws.Range("S23").Value = 31.8
$ws.Range("T23").Value = 41.8
$ws.Range("U23").Value = 20.6
$ws.Range("Y23").Value = 5.4
$ws.Range("AB23").Value = 95.7
$ws.Range("AC23").Value = -5.7
$ws.Range("AD23").Value = 1
$ws.Range("AJ23").Value = 20
$ws.Range("AL23").Value = 22
$ws.Range("AR23").Value = 25
$ws.Range("AY23").Value = 24
$ws.Range("AZ23").Value = 19
$ws.Range("AH24").Value = 20
$ws.Range("AJ24").Value = 22
$ws.Range("AH25").Value = 8
$ws.Range("AJ25").Value = 6
$ws.Range("AX26").Value = 18
$ws.Range("AW27").Value = 26
$ws.Range("D28").Value = 82
$ws.Range("E28").Value = 55
$ws.Range("G28").Value = 0.671
$ws.Range("J28").Value = 83.59999999999999
$ws.Range("K28").Value = 0.468
$ws.Range("L28").Value = 8.300000000000001
$ws.Range("N28").Value = 0.367
$ws.Range("P28").Value = 21.4
$ws.Range("S28").Value = 33.8
$ws.Range("T28").Value = 43.6
$ws.Range("V28").Value = 14
$ws.Range("AC28").Value = 6.2
$ws.Range("AD28").Value = 1
$ws.Range("AE28").Value = 5
$ws.Range("AG28").Value = 5
$ws.Range("AP28").Value = 23
$ws.Range("AT28").Value = 15
$ws.Range("AV28").Value = 13
$ws.Range("AZ29").Value = 18
$ws.Range("AP30").Value = 13
$ws.Range("AS30").Value = 20
$ws.Range("AT30").Value = 11
$ws.Range("AH31").Value = 1
$ws.Range("AJ31").Value = 21
$ws.Range("AP31").Value = 22
$ws.Range("BB31").Value = 17
$ws.Range("BC31").Value = 14
